# Add the 2023 season rows to the Money League historical-stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns: Manager, Finish, Wins, Points, PointsAllowed, Playoffs, DraftPosition, Champion
$newRows = @(
    @("Colin",    1, 13, 2157.94, 1686.66, 1, 11, 1),
    @("Mike",     2,  9, 1776.48, 1703.04, 1,  3, 0),
    @("Charles",  3,  8, 1812.62, 1787.1,  1,  9, 0),
    @("Alex",     4,  9, 1837.76, 1779.74, 1,  4, 1),
    @("Marcus",   5,  6, 1925.24, 1918.22, 1,  8, 0),
    @("EricR",    6,  7, 1823.16, 1873.2,  1, 12, 0),
    @("Chris",    7,  6, 1839.92, 1790.86, 0, 10, 0),
    @("Erik",     8,  6, 1842.76, 1886.36, 0,  7, 0),
    @("EricNC",   9,  6, 1831.6,  1829.2,  0,  6, 0),
    @("Jennifer",10,  5, 1710.6,  2018.88, 0,  2, 0),
    @("Chester", 11,  4, 1852.38, 2062.16, 0,  5, 0),
    @("John",    12,  5, 1723.32, 1798.36, 0,  1, 0)
)

$startRow = 138
$endRow = $startRow + $newRows.Count - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = 2023
    $ws.Cells.Item($row, 2).Value = $data[0]
    $ws.Cells.Item($row, 3).Value = $data[1]
    $ws.Cells.Item($row, 4).Value = $data[2]
    $ws.Cells.Item($row, 5).Value = $data[3]
    $ws.Cells.Item($row, 6).Value = $data[4]
    $ws.Cells.Item($row, 7).Value = $data[5]
    $ws.Cells.Item($row, 8).Value = $data[6]
    $ws.Cells.Item($row, 9).Value = $data[7]

    # Match the existing table's cell style (centered alignment).
    $rowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 9))
    $rowRange.HorizontalAlignment = -4108
}

# Fill the J-column formula down across all of the newly added rows in one
# shot, the same way Excel would if you dragged the fill handle from J137.
$formulaRange = $ws.Range($ws.Cells.Item($startRow, 10), $ws.Cells.Item($endRow, 10))
$formulaRange.Formula = "=IF(C" + $startRow + "=1,1,0)"
$formulaRange.HorizontalAlignment = -4108

$ws.Range("G149").Select()
